# Update LR-pair rows 2-5 (data rows 1-4) with the new TPM-based values.
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
#          E..T = numeric metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Cntn1 -> Ptprz1 -> FAPs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cntn1"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("I2").Value = 0.3500709860127268
$ws.Range("J2").Value = 0.446885632088942
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.01638633333333333
$ws.Range("N2").Value = 0.049159
$ws.Range("O2").Value = 0.02187172081577483
$ws.Range("P2").Value = 0.03245268321021395
$ws.Range("Q2").Value = 0.001152215952555556
$ws.Range("R2").Value = 0.010369943573
$ws.Range("S2").Value = 0.007656654871773377
$ws.Range("T2").Value = 0.01450263784937866

# Row 3: FAPs -> Cntn1 -> Ptprz1 -> MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cntn1"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("I3").Value = 0.3500709860127268
$ws.Range("J3").Value = 0.446885632088942
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7328155000000001
$ws.Range("N3").Value = 1.465631
$ws.Range("O3").Value = 0.9781282791842253
$ws.Range("P3").Value = 0.9675473167897861
$ws.Range("Q3").Value = 0.05152841042616667
$ws.Range("R3").Value = 0.309170462557
$ws.Range("S3").Value = 0.3424143311409535
$ws.Range("T3").Value = 0.4323829942395634

# Row 4: MuSCs -> Cntn1 -> Ptprz1 -> FAPs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Cntn1"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1305455
$ws.Range("H4").Value = 0.261091
$ws.Range("I4").Value = 0.6499290139872732
$ws.Range("J4").Value = 0.5531143679110581
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.01638633333333333
$ws.Range("N4").Value = 0.049159
$ws.Range("O4").Value = 0.02187172081577483
$ws.Range("P4").Value = 0.03245268321021395
$ws.Range("Q4").Value = 0.002139162078166667
$ws.Range("R4").Value = 0.012834972469
$ws.Range("S4").Value = 0.01421506594400145
$ws.Range("T4").Value = 0.0179500453608353

# Row 5: MuSCs -> Cntn1 -> Ptprz1 -> MuSCs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Cntn1"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1305455
$ws.Range("H5").Value = 0.261091
$ws.Range("I5").Value = 0.6499290139872732
$ws.Range("J5").Value = 0.5531143679110581
$ws.Range("M5").Value = 0.7328155000000001
$ws.Range("N5").Value = 1.465631
$ws.Range("O5").Value = 0.9781282791842253
$ws.Range("P5").Value = 0.9675473167897861
$ws.Range("Q5").Value = 0.09566576585525001
$ws.Range("R5").Value = 0.382663063421
$ws.Range("S5").Value = 0.6357139480432719
$ws.Range("T5").Value = 0.5351643225502228
